$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 1387.0333
$ws.Range("I40").Value = 1251.2
$ws.Range("K40").Value = 1251.2
$ws.Range("M40").Value = -1076.2

$ws.Range("H62").Value = 71451710
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 125038250
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 125038250
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -125039498

$ws.Range("H64").Value = 2881.9092
$ws.Range("I64").Value = 2816.6667
$ws.Range("J64").Value = 2960.2
$ws.Range("K64").Value = 2816.6667
$ws.Range("L64").Value = 2960.2
$ws.Range("M64").Value = -2568.6667
$ws.Range("N64").Value = -3456.2

$ws.Range("H65").Value = 71451710
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 125038250
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 625191250
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -625197490

$ws.Range("H67").Value = 2881.9092
$ws.Range("I67").Value = 2816.6667
$ws.Range("J67").Value = 2960.2
$ws.Range("K67").Value = 2816.6667
$ws.Range("L67").Value = 2960.2
$ws.Range("M67").Value = -1958.6667
$ws.Range("N67").Value = -4676.2

$ws.Range("H137").Value = 2420334.8
$ws.Range("I137").Value = 1112067.9
$ws.Range("K137").Value = 3336203.7
$ws.Range("M137").Value = -3333653.7

$ws.Range("H138").Value = 1474.27
$ws.Range("I138").Value = 470.71796
$ws.Range("J138").Value = 2115.8853
$ws.Range("K138").Value = 1412.15388
$ws.Range("L138").Value = 6347.6559
$ws.Range("M138").Value = 3727.84612
$ws.Range("N138").Value = -16627.6559

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H45").Value = 3781.4546
$ws.Range("I45").Value = 2599.923
$ws.Range("K45").Value = 2599.923
$ws.Range("M45").Value = -2222.923

$ws.Range("H61").Value = 1444.8889
$ws.Range("I61").Value = 1547.5385
$ws.Range("J61").Value = 1178
$ws.Range("K61").Value = 1547.5385
$ws.Range("L61").Value = 1178
$ws.Range("M61").Value = -1335.5385
$ws.Range("N61").Value = -1602

$ws.Range("H74").Value = 5129137.5
$ws.Range("I74").Value = 928.7027
$ws.Range("K74").Value = 928.7027
$ws.Range("M74").Value = -54.70270000000005

$ws.Range("H77").Value = 5129137.5
$ws.Range("I77").Value = 928.7027
$ws.Range("K77").Value = 4643.5135
$ws.Range("M77").Value = -275.5135

$ws.Range("H132").Value = 72163.64999999999
$ws.Range("I132").Value = 90481.44500000001
$ws.Range("K132").Value = 271444.335
$ws.Range("M132").Value = -268914.335

$ws.Range("H136").Value = 1444.8889
$ws.Range("I136").Value = 1547.5385
$ws.Range("J136").Value = 1178
$ws.Range("K136").Value = 4642.6155
$ws.Range("L136").Value = 3534
$ws.Range("M136").Value = -2092.6155
$ws.Range("N136").Value = -8634

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 46575.105
$ws.Range("I134").Value = 53599.668
$ws.Range("K134").Value = 160799.004
$ws.Range("M134").Value = -158264.004

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2687.9473
$ws.Range("I31").Value = 1906.9166
$ws.Range("K31").Value = 1906.9166
$ws.Range("M31").Value = -1611.9166

$ws.Range("H34").Value = 2687.9473
$ws.Range("I34").Value = 1906.9166
$ws.Range("K34").Value = 1906.9166
$ws.Range("M34").Value = -1704.9166

$ws.Range("H58").Value = 650.7432
$ws.Range("I58").Value = 695.8269
$ws.Range("J58").Value = 544.1818
$ws.Range("K58").Value = 695.8269
$ws.Range("L58").Value = 544.1818
$ws.Range("M58").Value = -492.8269
$ws.Range("N58").Value = -950.1818

$ws.Range("H62").Value = 3488.4614
$ws.Range("I62").Value = 3538.889
$ws.Range("J62").Value = 3375
$ws.Range("K62").Value = 3538.889
$ws.Range("L62").Value = 3375
$ws.Range("M62").Value = -2914.889
$ws.Range("N62").Value = -4623

$ws.Range("H65").Value = 3488.4614
$ws.Range("I65").Value = 3538.889
$ws.Range("J65").Value = 3375
$ws.Range("K65").Value = 17694.445
$ws.Range("L65").Value = 16875
$ws.Range("M65").Value = -14574.445
$ws.Range("N65").Value = -23115

$ws.Range("H132").Value = 1614.527
$ws.Range("I132").Value = 1454.9517
$ws.Range("J132").Value = 2439
$ws.Range("K132").Value = 4364.855100000001
$ws.Range("L132").Value = 7317
$ws.Range("M132").Value = -1834.855100000001
$ws.Range("N132").Value = -12377

$ws.Range("H134").Value = 2384.075
$ws.Range("I134").Value = 2479.257
$ws.Range("J134").Value = 1717.8
$ws.Range("K134").Value = 7437.771000000001
$ws.Range("L134").Value = 5153.4
$ws.Range("M134").Value = -4902.771000000001
$ws.Range("N134").Value = -10223.4

$ws.Range("H136").Value = 650.7432
$ws.Range("I136").Value = 695.8269
$ws.Range("J136").Value = 544.1818
$ws.Range("K136").Value = 2087.4807
$ws.Range("L136").Value = 1632.5454
$ws.Range("M136").Value = 462.5192999999999
$ws.Range("N136").Value = -6732.5454

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H107").Value = 1324.5454
$ws.Range("I107").Value = 701.6667
$ws.Range("J107").Value = 2072
$ws.Range("K107").Value = 2105.0001
$ws.Range("L107").Value = 6216
$ws.Range("M107").Value = -185.0001000000002
$ws.Range("N107").Value = -10056

$ws.Range("H131").Value = 897.96387
$ws.Range("I131").Value = 494.91666
$ws.Range("J131").Value = 966.08453
$ws.Range("K131").Value = 1484.74998
$ws.Range("L131").Value = 2898.25359
$ws.Range("M131").Value = 3555.25002
$ws.Range("N131").Value = -12978.25359

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H126").Value = 4406.162
$ws.Range("I126").Value = 2526.5557
$ws.Range("J126").Value = 6186.8423
$ws.Range("K126").Value = 7579.6671
$ws.Range("L126").Value = 18560.5269
$ws.Range("M126").Value = -5109.6671
$ws.Range("N126").Value = -23500.5269

$ws.Range("H132").Value = 1453.9647
$ws.Range("I132").Value = 1115.569
$ws.Range("K132").Value = 3346.707
$ws.Range("M132").Value = -816.7069999999999

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H136").Value = 1623.75
$ws.Range("I136").Value = 1510.5312
$ws.Range("J136").Value = 1850.1875
$ws.Range("K136").Value = 4531.5936
$ws.Range("L136").Value = 5550.5625
$ws.Range("M136").Value = -1981.5936
$ws.Range("N136").Value = -10650.5625

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H132").Value = 1258.5416
$ws.Range("I132").Value = 1083.8358
$ws.Range("J132").Value = 3599.6
$ws.Range("K132").Value = 3251.5074
$ws.Range("L132").Value = 10798.8
$ws.Range("M132").Value = -721.5074000000004
$ws.Range("N132").Value = -15858.8

$ws.Range("H136").Value = 1290.2295
$ws.Range("I136").Value = 1290.2295
$ws.Range("K136").Value = 3870.6885
$ws.Range("M136").Value = -1320.6885

